# feat: add 2022-Q4 data
#
# Before: sheet "总计" (summary) + sheet "2022-Q3" (fund detail for Q3).
# After:  sheet "总计" gets a new Q4 summary row (old Q3 row pushed down);
#         the "2022-Q3" sheet is renamed to "2022-Q4" and repopulated with
#         the Q4 fund detail; a fresh "2022-Q3" sheet is inserted after it
#         holding the original Q3 fund detail that used to live there.

$wb = $excel.ActiveWorkbook

$zong = $wb.Worksheets.Item("总计")
$oldQ3 = $wb.Worksheets.Item("2022-Q3")

# --- 1. Make room in 总计: push the existing Q3 summary row down to row 3,
#        write the new Q4 summary numbers into row 2. -----------------------
$zong.Range("A2").Copy()
$zong.Range("A3").PasteSpecial(-4122)

$zong.Range("A3").Value = 1
$zong.Range("B3").Value = "2022-Q3"
$zong.Range("C3").Value = 2
$zong.Range("D3").Value = 0.32

$zong.Range("B2").Value = "2022-Q4"
$zong.Range("C2").Value = 4
$zong.Range("D2").Value = 0.93

# --- 2. Turn the old "2022-Q3" sheet into the new "2022-Q4" sheet, but
#        first clone its current (Q3) content into a brand-new sheet so the
#        Q3 detail survives under its own tab. ------------------------------
$oldQ3.Name = "2022-Q4"
$q4 = $wb.Worksheets.Item("2022-Q4")

$newQ3 = $wb.Worksheets.Add($null, $q4)
$newQ3.Name = "2022-Q3"
$q4.Range("A1:H3").Copy($newQ3.Range("A1"))
$newQ3.Range("A1").ClearContents()

# --- 3. Overwrite the (now) "2022-Q4" sheet with the Q4 fund data, reusing
#        总计's header style (style index 2) for the label/index cells. -----
$zong.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$q4.Range("A2:A5").PasteSpecial(-4122)

function Set-TextCell($range, $text) {
    # Route the literal through a text formula, then bake it back down to a
    # plain value - this keeps leading zeros / decimal-looking strings as
    # genuine text (no silent "012584" -> 12584 numeric coercion) without
    # picking up a stray NumberFormat/quote-prefix style along the way.
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

$q4.Range("A2").Value = 0
Set-TextCell $q4.Range("B2") "202801"
Set-TextCell $q4.Range("C2") "南方全球精选配置（QDII-FOF）"
Set-TextCell $q4.Range("D2") "17.02"
Set-TextCell $q4.Range("E2") "32.64"
Set-TextCell $q4.Range("F2") "3.51"
Set-TextCell $q4.Range("G2") "0.5974"
$q4.Range("H2").Value = 1

$q4.Range("A3").Value = 1
Set-TextCell $q4.Range("B3") "160125"
Set-TextCell $q4.Range("C3") "南方香港优选股票（QDII-LOF）"
Set-TextCell $q4.Range("D3") "2.27"
Set-TextCell $q4.Range("E3") "84.75"
Set-TextCell $q4.Range("F3") "7.42"
Set-TextCell $q4.Range("G3") "0.1684"
$q4.Range("H3").Value = 1

$q4.Range("A4").Value = 2
Set-TextCell $q4.Range("B4") "012584"
Set-TextCell $q4.Range("C4") "南方中国新兴经济9个月持有期混合（QDII）A"
Set-TextCell $q4.Range("D4") "3.08"
Set-TextCell $q4.Range("E4") "83.21"
Set-TextCell $q4.Range("F4") "5.26"
Set-TextCell $q4.Range("G4") "0.1620"
$q4.Range("H4").Value = 3

$q4.Range("A5").Value = 3
Set-TextCell $q4.Range("B5") "012585"
Set-TextCell $q4.Range("C5") "南方中国新兴经济9个月持有期混合（QDII）C"
Set-TextCell $q4.Range("D5") "0.12"
Set-TextCell $q4.Range("E5") "83.21"
Set-TextCell $q4.Range("F5") "5.26"
Set-TextCell $q4.Range("G5") "0.0063"
$q4.Range("H5").Value = 3
